$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 37.00997885040392
$ws.Range("C2").Value = 38.29388288199763
$ws.Range("D2").Value = 35.67999920255267
$ws.Range("E2").Value = 37.01077875011296
$ws.Range("F2").Value = 37.15300573523157
$ws.Range("G2").Value = 36.54477116405194
$ws.Range("H2").Value = 40.05944957231277
$ws.Range("I2").Value = 32.59176978264608
$ws.Range("J2").Value = 37.12425180097848
$ws.Range("K2").Value = 36.92397212398041
$ws.Range("L2").Value = 37.19706388087692
$ws.Range("M2").Value = 36.34222408766744
$ws.Range("N2").Value = 17.43473966858451
$ws.Range("O2").Value = 31.91991872074886
$ws.Range("P2").Value = 40.9707364202819
$ws.Range("Q2").Value = 33.90553370188189
